$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newDate = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$newDate = $newDate.AddDays(45224)

$ws.Range("C2").Value = $newDate
$ws.Range("C3").Value = $newDate
$ws.Range("C4").Value = $newDate
$ws.Range("C5").Value = $newDate
